$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (rich-text shared strings) ---
$ws.Range("A8").Characters(21,1).Text = "2"
$ws.Range("C9").Characters(27,8).Text = "1/8/2024"
$ws.Range("C9").Characters(46,8).Text = "1/14/2024"

# --- Data table updates ---
$ws.Range("G15").Value = 1
$ws.Range("N15").Value = -100
$ws.Range("K37").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C28").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("E28").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 2
$ws.Range("K16").Value = -33.333333333333
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -75
$ws.Range("N16").Value = -95.121951219512
$ws.Range("C17").Value = 1
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0"
$ws.Range("C28").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "***.*"
$ws.Range("E28").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 200
$ws.Range("L17").Value = -33.333333333333
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -15.384615384615
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 9
$ws.Range("K18").Value = -22.222222222222
$ws.Range("L18").Value = 75
$ws.Range("M18").Value = 16.666666666666
$ws.Range("N18").Value = -88.135593220339
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -13.043478260869
$ws.Range("I19").Value = 11
$ws.Range("J19").Value = 10
$ws.Range("K19").Value = 10
$ws.Range("L19").Value = -47.619047619047
$ws.Range("M19").Value = -38.888888888888
$ws.Range("N19").Value = -60.714285714285
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -80
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = -75
$ws.Range("N20").Value = -98.907103825136
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -38.888888888888
$ws.Range("F21").Value = 59
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = 1.724137931034
$ws.Range("I21").Value = 26
$ws.Range("J21").Value = 31
$ws.Range("K21").Value = -16.129032258064
$ws.Range("L21").Value = -23.529411764705
$ws.Range("M21").Value = -27.777777777777
$ws.Range("N21").Value = -91.666666666666
$ws.Range("C22").Value = 1
$ws.Range("C36").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("C36").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = 0
$ws.Range("K37").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("I22").Value = 1
$ws.Range("C36").Copy()
$ws.Range("I22").PasteSpecial(-4122)
$ws.Range("J22").Value = 1
$ws.Range("C36").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("K22").Value = 0
$ws.Range("K37").Copy()
$ws.Range("K22").PasteSpecial(-4122)
$ws.Range("L22").Value = -66.666666666666
$ws.Range("M22").Value = -50
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C28").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 7.407407407407
$ws.Range("F24").Value = 124
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = 16.981132075471
$ws.Range("I24").Value = 52
$ws.Range("J24").Value = 49
$ws.Range("K24").Value = 6.122448979591
$ws.Range("L24").Value = -5.454545454545
$ws.Range("M24").Value = 100
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 23
$ws.Range("H25").Value = 9.523809523809
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 10
$ws.Range("K25").Value = -50
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -37.5
$ws.Range("G26").Value = 1
$ws.Range("L26").Value = -100
$ws.Range("K37").Copy()
$ws.Range("L26").PasteSpecial(-4122)
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C28").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("C36").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("K37").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("C36").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H27").Value = 200
$ws.Range("K37").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("J27").Value = 1
$ws.Range("C36").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("K27").Value = 100
$ws.Range("K37").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("J41").Value = 475
$ws.Range("K41").Value = -28.892215568862
$ws.Range("L41").Value = -35.897435897435
$ws.Range("M41").Value = -53.793774319066
$ws.Range("N41").Value = -64.55223880597
$ws.Range("J43").Value = 953
$ws.Range("K43").Value = -51.228249744114
$ws.Range("L43").Value = -63.846737481031
$ws.Range("M43").Value = -85.562793516133
$ws.Range("N43").Value = -87.030484485574

$excel.CutCopyMode = $false
